# Remove the "2 Quantil (Mediana)" example slide (sldId 279) from the deck.
# This was slide #12 - a worked example showing quantile calculations over
# a small data set ("Ejemplo:" / "Datos: 6, 7, 15, 36, ...") together with
# four supporting text boxes and click-triggered entrance animations.
$p = $ppt.ActivePresentation
$p.Slides.Item(12).Delete()
